# Regenerate the "three-digit ÷ one-digit" answer key: each table cell's
# division problem is replaced with a newly generated one. Every old
# problem string is unique in the document, so a plain Find/Replace
# (MatchWholeWord, not a wildcard pattern) per cell is safe and targeted.
$d = $word.ActiveDocument

$d.Content.Find.Execute("275÷6=45, 5", $true, $false, $false, $false, $false, $true, 1, $false, "124÷4=31, 0", 2) | Out-Null
$d.Content.Find.Execute("446÷3=148, 2", $true, $false, $false, $false, $false, $true, 1, $false, "740÷8=92, 4", 2) | Out-Null
$d.Content.Find.Execute("720÷4=180, 0", $true, $false, $false, $false, $false, $true, 1, $false, "876÷5=175, 1", 2) | Out-Null
$d.Content.Find.Execute("532÷7=76, 0", $true, $false, $false, $false, $false, $true, 1, $false, "501÷5=100, 1", 2) | Out-Null
$d.Content.Find.Execute("836÷6=139, 2", $true, $false, $false, $false, $false, $true, 1, $false, "406÷4=101, 2", 2) | Out-Null
$d.Content.Find.Execute("154÷2=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "764÷5=152, 4", 2) | Out-Null
$d.Content.Find.Execute("353÷8=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "259÷7=37, 0", 2) | Out-Null
$d.Content.Find.Execute("941÷6=156, 5", $true, $false, $false, $false, $false, $true, 1, $false, "525÷7=75, 0", 2) | Out-Null
$d.Content.Find.Execute("342÷9=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "202÷6=33, 4", 2) | Out-Null
$d.Content.Find.Execute("959÷2=479, 1", $true, $false, $false, $false, $false, $true, 1, $false, "754÷3=251, 1", 2) | Out-Null
$d.Content.Find.Execute("337÷7=48, 1", $true, $false, $false, $false, $false, $true, 1, $false, "420÷6=70, 0", 2) | Out-Null
$d.Content.Find.Execute("407÷5=81, 2", $true, $false, $false, $false, $false, $true, 1, $false, "151÷3=50, 1", 2) | Out-Null
$d.Content.Find.Execute("754÷4=188, 2", $true, $false, $false, $false, $false, $true, 1, $false, "960÷2=480, 0", 2) | Out-Null
$d.Content.Find.Execute("716÷8=89, 4", $true, $false, $false, $false, $false, $true, 1, $false, "565÷7=80, 5", 2) | Out-Null
$d.Content.Find.Execute("624÷7=89, 1", $true, $false, $false, $false, $false, $true, 1, $false, "188÷9=20, 8", 2) | Out-Null
$d.Content.Find.Execute("797÷2=398, 1", $true, $false, $false, $false, $false, $true, 1, $false, "532÷9=59, 1", 2) | Out-Null
$d.Content.Find.Execute("518÷9=57, 5", $true, $false, $false, $false, $false, $true, 1, $false, "208÷4=52, 0", 2) | Out-Null
$d.Content.Find.Execute("116÷6=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "402÷8=50, 2", 2) | Out-Null
$d.Content.Find.Execute("115÷2=57, 1", $true, $false, $false, $false, $false, $true, 1, $false, "908÷2=454, 0", 2) | Out-Null
$d.Content.Find.Execute("594÷4=148, 2", $true, $false, $false, $false, $false, $true, 1, $false, "965÷2=482, 1", 2) | Out-Null
$d.Content.Find.Execute("572÷2=286, 0", $true, $false, $false, $false, $false, $true, 1, $false, "550÷3=183, 1", 2) | Out-Null
$d.Content.Find.Execute("871÷8=108, 7", $true, $false, $false, $false, $false, $true, 1, $false, "946÷7=135, 1", 2) | Out-Null
$d.Content.Find.Execute("833÷3=277, 2", $true, $false, $false, $false, $false, $true, 1, $false, "118÷9=13, 1", 2) | Out-Null
$d.Content.Find.Execute("838÷8=104, 6", $true, $false, $false, $false, $false, $true, 1, $false, "637÷8=79, 5", 2) | Out-Null
$d.Content.Find.Execute("825÷3=275, 0", $true, $false, $false, $false, $false, $true, 1, $false, "738÷7=105, 3", 2) | Out-Null
